$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in Cameron's bio and extend it with a new closing sentence.
$ws.Range("B13").Value = "I finished my undergraduate degree at Laurentian University studying deer mice in Algonquin Provincial Park. Now I'm excited to make the move to aquatic ecosystems and to learn everything I can in this field! I am currently working on biodiversity in benthic stream ecosystems."

# Add a new "Educational Background" row for Cameron's current MSc studies.
$ws.Range("A21").Value = "MSc Biology"
$ws.Range("B21").Value = "Laurentian University"
$ws.Range("C21").Value = "current"

# Match the author's final active selection.
$ws.Range("C21").Select()
